$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Remove the whole row 28 (INDICATOR_59 entry) - this shifts all rows below it up by one,
# matching the target layout where the INDICATOR_59 row no longer exists.
$ws.Rows.Item(28).Delete()

# Restore the selection/view state to match the target workbook.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A28:XFD28").Select()
